$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 8735.166999999999
$ws.Range("J40").Value = 3323.5
$ws.Range("L40").Value = 3323.5
$ws.Range("N40").Value = -3673.5
$ws.Range("H86").Value = 131946640
$ws.Range("I86").Value = 200001970
$ws.Range("K86").Value = 200001970
$ws.Range("M86").Value = -200000847
$ws.Range("H89").Value = 131946640
$ws.Range("I89").Value = 200001970
$ws.Range("K89").Value = 1000009850
$ws.Range("M89").Value = -1000004234
$ws.Range("H112").Value = 14066.866
$ws.Range("J112").Value = 14066.866
$ws.Range("L112").Value = 42200.598
$ws.Range("N112").Value = -44416.598
$ws.Range("H135").Value = 385978.97
$ws.Range("I135").Value = 401018.16
$ws.Range("K135").Value = 3609163.44
$ws.Range("M135").Value = -3606628.44
$ws.Range("H138").Value = 1643681.1
$ws.Range("J138").Value = 2330249.2
$ws.Range("L138").Value = 6990747.600000001
$ws.Range("N138").Value = -7001027.600000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1869595.2
$ws.Range("I32").Value = 2196371.8
$ws.Range("K32").Value = 2196371.8
$ws.Range("M32").Value = -2196084.8
$ws.Range("H61").Value = 6694.702
$ws.Range("I61").Value = 3623.6667
$ws.Range("J61").Value = 10840.6
$ws.Range("K61").Value = 3623.6667
$ws.Range("L61").Value = 10840.6
$ws.Range("M61").Value = -3411.6667
$ws.Range("N61").Value = -11264.6
$ws.Range("H74").Value = 37164.5
$ws.Range("I74").Value = 47406.137
$ws.Range("K74").Value = 47406.137
$ws.Range("M74").Value = -46532.137
$ws.Range("H77").Value = 37164.5
$ws.Range("I77").Value = 47406.137
$ws.Range("K77").Value = 237030.685
$ws.Range("M77").Value = -232662.685
$ws.Range("H132").Value = 1511274.1
$ws.Range("I132").Value = 3515903.8
$ws.Range("K132").Value = 10547711.4
$ws.Range("M132").Value = -10545181.4
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 6694.702
$ws.Range("I136").Value = 3623.6667
$ws.Range("J136").Value = 10840.6
$ws.Range("K136").Value = 10871.0001
$ws.Range("L136").Value = 32521.8
$ws.Range("M136").Value = -8321.000100000001
$ws.Range("N136").Value = -37621.8

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 69330
$ws.Range("J76").Value = 69330
$ws.Range("L76").Value = 69330
$ws.Range("N76").Value = -69960
$ws.Range("H79").Value = 69330
$ws.Range("J79").Value = 69330
$ws.Range("L79").Value = 69330
$ws.Range("N79").Value = -71514
$ws.Range("H105").Value = 3833.3
$ws.Range("J105").Value = 5765.25
$ws.Range("L105").Value = 5765.25
$ws.Range("N105").Value = -9259.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6160.5146
$ws.Range("I31").Value = 3491.7837
$ws.Range("K31").Value = 3491.7837
$ws.Range("M31").Value = -3196.7837
$ws.Range("H34").Value = 6160.5146
$ws.Range("I34").Value = 3491.7837
$ws.Range("K34").Value = 3491.7837
$ws.Range("M34").Value = -3289.7837
$ws.Range("H51").Value = 63984.332
$ws.Range("J51").Value = 63984.332
$ws.Range("L51").Value = 63984.332
$ws.Range("N51").Value = -65456.332
$ws.Range("H61").Value = 63984.332
$ws.Range("J61").Value = 63984.332
$ws.Range("L61").Value = 63984.332
$ws.Range("N61").Value = -64680.332
$ws.Range("H134").Value = 4312.017
$ws.Range("I134").Value = 2567.5112
$ws.Range("J134").Value = 9919.357
$ws.Range("K134").Value = 7702.5336
$ws.Range("L134").Value = 29758.071
$ws.Range("M134").Value = -5167.5336
$ws.Range("N134").Value = -34828.071

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2323.973
$ws.Range("J68").Value = 2402.8386
$ws.Range("L68").Value = 7208.5158
$ws.Range("N68").Value = -8830.515800000001
$ws.Range("H71").Value = 2323.973
$ws.Range("J71").Value = 2402.8386
$ws.Range("L71").Value = 21625.5474
$ws.Range("N71").Value = -29737.5474
$ws.Range("H75").Value = 1350.7693
$ws.Range("I75").Value = 1251.5
$ws.Range("J75").Value = 1368.8182
$ws.Range("K75").Value = 3754.5
$ws.Range("L75").Value = 4106.4546
$ws.Range("M75").Value = -2756.5
$ws.Range("N75").Value = -6102.4546
$ws.Range("H78").Value = 1350.7693
$ws.Range("I78").Value = 1251.5
$ws.Range("J78").Value = 1368.8182
$ws.Range("K78").Value = 11263.5
$ws.Range("L78").Value = 12319.3638
$ws.Range("M78").Value = -6271.5
$ws.Range("N78").Value = -22303.3638
$ws.Range("H107").Value = 7409195
$ws.Range("J107").Value = 9525954
$ws.Range("L107").Value = 28577862
$ws.Range("N107").Value = -28581702
$ws.Range("H129").Value = 63940.938
$ws.Range("I129").Value = 974.875
$ws.Range("J129").Value = 126907
$ws.Range("K129").Value = 2924.625
$ws.Range("L129").Value = 380721
$ws.Range("M129").Value = 2075.375
$ws.Range("N129").Value = -390721

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6919.2666
$ws.Range("I70").Value = 4983.3335
$ws.Range("K70").Value = 4983.3335
$ws.Range("M70").Value = -4713.3335
$ws.Range("H73").Value = 6919.2666
$ws.Range("I73").Value = 4983.3335
$ws.Range("K73").Value = 4983.3335
$ws.Range("M73").Value = -4047.3335
$ws.Range("H122").Value = 83422830
$ws.Range("I122").Value = 125130250
$ws.Range("K122").Value = 375390750
$ws.Range("M122").Value = -375388300
$ws.Range("H132").Value = 3828.9429
$ws.Range("I132").Value = 2440.1738
$ws.Range("J132").Value = 6490.75
$ws.Range("K132").Value = 7320.5214
$ws.Range("L132").Value = 19472.25
$ws.Range("M132").Value = -4790.5214
$ws.Range("N132").Value = -24532.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2009.8889
$ws.Range("I16").Value = 1691
$ws.Range("J16").Value = 2647.6667
$ws.Range("K16").Value = 1691
$ws.Range("L16").Value = 2647.6667
$ws.Range("M16").Value = -1521
$ws.Range("N16").Value = -2987.6667
$ws.Range("H40").Value = 5168.6665
$ws.Range("I40").Value = 5146
$ws.Range("J40").Value = 5452
$ws.Range("K40").Value = 5146
$ws.Range("L40").Value = 5452
$ws.Range("M40").Value = -5010
$ws.Range("N40").Value = -5724
$ws.Range("H68").Value = 2894.5
$ws.Range("I68").Value = 2789
$ws.Range("K68").Value = 2789
$ws.Range("M68").Value = -2040
$ws.Range("H71").Value = 2894.5
$ws.Range("I71").Value = 2789
$ws.Range("K71").Value = 13945
$ws.Range("M71").Value = -10201
$ws.Range("H132").Value = 12826658
$ws.Range("I132").Value = 22730308
$ws.Range("J132").Value = 10170.177
$ws.Range("K132").Value = 68190924
$ws.Range("L132").Value = 30510.531
$ws.Range("M132").Value = -68188394
$ws.Range("N132").Value = -35570.531

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 623
$ws.Range("I100").Value = 541.25
$ws.Range("J100").Value = 950
$ws.Range("K100").Value = 1082.5
$ws.Range("L100").Value = 1900
$ws.Range("M100").Value = -541.5
$ws.Range("N100").Value = -2982
$ws.Range("H122").Value = 3859.7058
$ws.Range("I122").Value = 3708.1052
$ws.Range("J122").Value = 4051.7334
$ws.Range("K122").Value = 11124.3156
$ws.Range("L122").Value = 12155.2002
$ws.Range("M122").Value = -8674.3156
$ws.Range("N122").Value = -17055.2002
$ws.Range("H126").Value = 3843
$ws.Range("I126").Value = 1125
$ws.Range("K126").Value = 3375
$ws.Range("M126").Value = -905
$ws.Range("H132").Value = 17870816
$ws.Range("I132").Value = 29415944
$ws.Range("K132").Value = 88247832
$ws.Range("M132").Value = -88245302
$ws.Range("H133").Value = 144975
$ws.Range("J133").Value = 144975
$ws.Range("L133").Value = 144975
$ws.Range("N133").Value = -155095
$ws.Range("H136").Value = 27030518
$ws.Range("I136").Value = 45456316
$ws.Range("J136").Value = 6013.2
$ws.Range("K136").Value = 136368948
$ws.Range("L136").Value = 18039.6
$ws.Range("M136").Value = -136366398
$ws.Range("N136").Value = -23139.6
$ws.Range("H139").Value = 85413.875
$ws.Range("J139").Value = 85413.875
$ws.Range("L139").Value = 85413.875
$ws.Range("N139").Value = -95693.875
